{"js": "// The commit changes the heading \"3. Facebook and privacy\" to\n// \"3. Social networks and privacy\" (i.e. \"Facebook\" -> \"Social networks\").\n// All other hunks in the source diff are inert re-serialization noise\n// (proofErr markers removed / runs re-split by Word) that do not change\n// the document's visible text, so only this replacement is applied here.\n\nconst body = context.document.body;\nconst results = body.search(\"Facebook\", { matchCase: true, matchWholeWord: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find 'Facebook' in the document body.\");\n}\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"Social networks\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# The commit changes the heading \"3. Facebook and privacy\" to\n# \"3. Social networks and privacy\" (i.e. \"Facebook\" -> \"Social networks\").\n# All other hunks in the source diff are inert re-serialization noise\n# (proofErr markers removed / runs re-split by Word) that do not change\n# the document's visible text, so only this replacement is applied here.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$found = $find.Execute(\n    \"Facebook\",         # FindText\n    $true,              # MatchCase\n    $true,              # MatchWholeWord\n    $false,             # MatchWildcards\n    $false,             # MatchSoundsLike\n    $false,             # MatchAllWordForms\n    $true,              # Forward\n    1,                  # Wrap (wdFindContinue)\n    $false,             # Format\n    \"Social networks\",  # ReplaceWith\n    2                   # Replace (wdReplaceAll)\n)\n\nif (-not $found) {\n    throw \"Could not find 'Facebook' in the document body.\"\n}\n"}
